# "fixed styles of basics"
# - Removes the hard-coded dark-blue (002060) direct character formatting that
#   was applied on top of the Heading3/Heading4 paragraphs (both on the
#   paragraph mark and on the run(s) of text), letting the headings inherit
#   their color from the style instead.
# - Splits the "La Pâte brisée" run into two runs ("La Pâte " + "brisée").
# - Updates the Heading 3 / Heading 4 (and their linked character styles)
#   definitions so the style itself now carries color 002060 (plain, not
#   theme-based) instead of C00000 / E36C0A(accent6).

$d = $word.ActiveDocument
$wNs = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'

# ---------------------------------------------------------------------------
# 1. document.xml - strip the direct w:color formatting from the heading
#    paragraphs and split the title run in two.
# ---------------------------------------------------------------------------

for ($i = $d.Paragraphs.Count; $i -ge 1; $i--) {
    $p = $d.Paragraphs($i)
    $styleName = $p.Style.NameLocal
    $text = $p.Range.Text.TrimEnd("`r")

    if ($styleName -eq "Heading 3" -and $text -eq "La Pâte brisée") {
        $frag = '<w:p ' + $wNs + '>' +
                  '<w:pPr><w:pStyle w:val="Heading3"/></w:pPr>' +
                  '<w:bookmarkStart w:id="0" w:name="_Toc395353023"/>' +
                  '<w:r><w:t xml:space="preserve">La Pâte </w:t></w:r>' +
                  '<w:r><w:t>brisée</w:t></w:r>' +
                  '<w:bookmarkEnd w:id="0"/>' +
                '</w:p>'
        $null = $p.Range.InsertXML($frag)
    }
    elseif ($styleName -eq "Heading 4" -and $text -eq "Ingrédients") {
        $frag = '<w:p ' + $wNs + '>' +
                  '<w:pPr><w:pStyle w:val="Heading4"/></w:pPr>' +
                  '<w:r><w:t>Ingrédients</w:t></w:r>' +
                '</w:p>'
        $null = $p.Range.InsertXML($frag)
    }
    elseif ($styleName -eq "Heading 4" -and $text -eq "Préparation") {
        $frag = '<w:p ' + $wNs + '>' +
                  '<w:pPr><w:pStyle w:val="Heading4"/></w:pPr>' +
                  '<w:r><w:t>Préparation</w:t></w:r>' +
                '</w:p>'
        $null = $p.Range.InsertXML($frag)
    }
    elseif ($styleName -eq "Heading 4" -and $text -eq "Remarque") {
        $frag = '<w:p ' + $wNs + '>' +
                  '<w:pPr><w:pStyle w:val="Heading4"/></w:pPr>' +
                  '<w:r><w:t>Remarque</w:t></w:r>' +
                '</w:p>'
        $null = $p.Range.InsertXML($frag)
    }
}

# ---------------------------------------------------------------------------
# 2. styles.xml - recolor Heading 3 / Heading 4 (and linked char styles)
#    from their old colors to plain 002060.
# ---------------------------------------------------------------------------

$newColor = 6299648   # BGR-packed value Word uses internally for RGB(0x00,0x20,0x60)

foreach ($styleName in @("Heading 3", "Heading 4", "Heading3Char", "Heading4Char")) {
    $s = $d.Styles($styleName)
    $s.Font.Color = $newColor
}

Write-Host "done"
